$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Press Release")

$rows = @(
  @{
    r=6
    A="SEBI"; B="Press Release"; C=2025; D="November"; E="2025-11-14"
    F="Inauguration of SEBI’s Pavilion at the IITF, 2025 by Chairman, SEBI"
    G="https://www.sebi.gov.in/sebi_data/attachdocs/nov-2025/1763122672135.pdf"
    H="1763122672135.pdf"
    I="/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Press Release/2025/November/1763122672135.pdf"
    J="No.74/2025"
    K="`"SEBI, in collaboration with market associations, has set up a pavilion named 'Bharat Ka Share Bazaar' at the 44th India International Trade Fair 2025 in New Delhi. The pavilion aims to promote investor protection, financial education, and market integrity through live demonstrations, expert talks, and interactive activities, with a focus on enhancing national financial literacy.`""
    L="The Securities and Exchange Board of India (SEBI), along with various market associations, has established a pavilion named `"Bharat Ka Share Bazaar`" at the 44th India International Trade Fair (IITF) 2025, held from November 14–27, 2025, in New Delhi. The pavilion focuses on investor protection, financial education, and market integrity, featuring live demonstrations, expert talks, and interactive activities to promote financial literacy nationwide."
  },
  @{
    r=7
    A="SEBI"; B="Press Release"; C=2025; D="November"; E="2025-11-11"
    F="SEBI and IEPFA Conducted the third “Niveshak Shivir” at Amritsar - Proactive  Joint Initiative to Minimize Unclaimed Assets and Facilitate Investor Claims for  Unpaid Dividends and Shares"
    G="https://www.sebi.gov.in/sebi_data/attachdocs/nov-2025/1762858201080.pdf"
    H="1762858201080.pdf"
    I="/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Press Release/2025/November/1762858201080.pdf"
    J="No.72/2025"
    K="`"SEBI, IEPFA, and Market Infrastructure Institutions held a three-day investor event in Amritsar on Nov 1st, 2025. The focus was on assisting investors to reclaim unpaid dividends and unclaimed shares. Service desks were established for claim assistance, KYC updates, and issue resolution. A collaborative video on claiming unpaid dividends was launched. Similar events are planned for cities with substantial investor assets, including Mumbai.`""
    L="The third `"Niveshak Shivir`" was jointly organized by SEBI, IEPFA, and key Market Infrastructure Institutions on November 01, 2025, in Amritsar, aiming to help investors reclaim unpaid dividends and unclaimed shares. Dedicated service desks were set up for claim assistance, KYC updates, and issue resolution. A collaborative explainer video on claiming unpaid dividends was also launched. Similar events are proposed for cities with significant unclaimed investor assets, such as Mumbai."
  },
  @{
    r=8
    A="SEBI"; B="Press Release"; C=2025; D="November"; E="2025-11-14"
    F="Inauguration of SEBI’s Pavilion at the IITF, 2025 by Chairman, SEBI"
    G="https://www.sebi.gov.in/sebi_data/attachdocs/nov-2025/1763122672135.pdf"
    H="1763122672135.pdf"
    I="/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Press Release/2025/November/1763122672135.pdf"
    J="No.74/2025"
    K="`"SEBI has launched 'Bharat Ka Share Bazaar' pavilion at the 44th IITF 2025, offering financial education and investor awareness sessions from November 14–27, 2025. This initiative is part of SEBI's ongoing commitment to market integrity and investor protection.`""
    L="The Securities and Exchange Board of India (SEBI), in collaboration with various financial entities, has established a pavilion titled `"Bharat Ka Share Bazaar`" at the 44th India International Trade Fair (IITF) 2025 from November 14–27, 2025. The pavilion aims to showcase financial education and investor awareness, with interactive demonstrations and services offered. This initiative is part of SEBI's continued efforts towards market integrity and investor protection."
  },
  @{
    r=9
    A="SEBI"; B="Press Release"; C=2025; D="November"; E="2025-11-12"
    F="Report of the High Level Committee on Conflict of Interest, Disclosures and  related matters in respect of Members and Officials of SEBI"
    G="https://www.sebi.gov.in/sebi_data/attachdocs/nov-2025/1762943978568.pdf"
    H="1762943978568.pdf"
    I="/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Press Release/2025/November/1762943978568.pdf"
    J="No. 73/2025"
    K="`"The High Level Committee has submitted a report regarding conflict of interest and disclosures for SEBI Members and Officials, available on the SEBI website as of November 10, 2025. This report aims to enhance transparency and accountability within SEBI.`""
    L="The High Level Committee has submitted a report on conflict of interest, disclosures, and related matters concerning Members and Officials of SEBI, dated November 10, 2025. The report is accessible on the SEBI website."
  },
  @{
    r=10
    A="SEBI"; B="Press Release"; C=2025; D="November"; E="2025-11-11"
    F="SEBI and IEPFA Conducted the third “Niveshak Shivir” at Amritsar - Proactive  Joint Initiative to Minimize Unclaimed Assets and Facilitate Investor Claims for  Unpaid Dividends and Shares"
    G="https://www.sebi.gov.in/sebi_data/attachdocs/nov-2025/1762858201080.pdf"
    H="1762858201080.pdf"
    I="/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Press Release/2025/November/1762858201080.pdf"
    J="No.72/2025"
    K="`"SEBI, IEPFA, and Market Infrastructure Institutions (MIIs) hosted a Niveshak Shivir in Amritsar on November 1, 2025. The event aims to help investors reclaim unpaid dividends and unclaimed shares, reduce the volume of unclaimed assets, and offer assistance for claim-related issues. Similar events are planned in other cities with significant unclaimed investor assets, including Mumbai.`""
    L="The third `"Niveshak Shivir`" was jointly organized by SEBI, IEPFA, and key Market Infrastructure Institutions (MIIs) on November 1, 2025, in Amritsar. The event aimed to help investors reclaim unpaid dividends and unclaimed shares, reduce the volume of unclaimed assets, and provide end-to-end assistance for claim-related issues. A series of similar events are proposed in other cities with significant unclaimed investor assets, including Mumbai."
  },
  @{
    r=11
    A="SEBI"; B="Press Release"; C=2025; D="November"; E="2025-11-10"
    F="High Level Committee on “Conflict of Interest, Disclosures and related matters in respect of Members and Officials of SEBI” submits its report"
    G="https://www.sebi.gov.in/sebi_data/attachdocs/nov-2025/1762786201875.pdf"
    H="1762786201875.pdf"
    I="/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Press Release/2025/November/1762786201875.pdf"
    J="No.71/2025"
    K="`"The Securities and Exchange Board of India (SEBI) has proposed updates to its policies concerning conflicts of interest, disclosures, and related matters for its Members and Officials. The changes are based on recommendations from a High Level Committee consisting of six members, including retired IAS officers, bank directors, and an ex-Auditor General. The aim is to establish a strong framework for managing potential conflicts, improve disclosure requirements, and address public concerns.`""
    L="The Securities and Exchange Board of India (SEBI) formed a High Level Committee (HLC) on March 24, 2025, to review and update policies related to conflict of interest, disclosures, and associated matters for its Members and Officials. The committee consisted of six members, including retired IAS officers, bank directors, and an ex-Auditor General. The committee submitted its report on November 10, 2025, recommending a robust framework for managing conflict of interest, disclosure requirements, and mechanisms to address public concerns."
  }
)

foreach ($row in $rows) {
  $r = $row.r
  # Column E holds plain-text date-like strings ("2025-11-14"); force text
  # formatting first so Excel doesn't auto-convert them into date serials.
  $ws.Cells.Item($r, 5).NumberFormat = "@"

  $ws.Cells.Item($r, 1).Value = $row.A
  $ws.Cells.Item($r, 2).Value = $row.B
  $ws.Cells.Item($r, 3).Value = $row.C
  $ws.Cells.Item($r, 4).Value = $row.D
  $ws.Cells.Item($r, 5).Value = $row.E
  $ws.Cells.Item($r, 6).Value = $row.F
  $ws.Cells.Item($r, 7).Value = $row.G
  $ws.Cells.Item($r, 8).Value = $row.H
  $ws.Cells.Item($r, 9).Value = $row.I
  $ws.Cells.Item($r, 10).Value = $row.J
  $ws.Cells.Item($r, 11).Value = $row.K
  $ws.Cells.Item($r, 12).Value = $row.L
}
